$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 121.05556
$ws.Range("I2").Value = 126.46667
$ws.Range("K2").Value = 126.46667
$ws.Range("M2").Value = -13.46666999999999
$ws.Range("H29").Value = 3481.3572
$ws.Range("J29").Value = 3752.875
$ws.Range("L29").Value = 11258.625
$ws.Range("N29").Value = -11820.625
$ws.Range("H55").Value = 898.5
$ws.Range("I55").Value = 898
$ws.Range("K55").Value = 898
$ws.Range("M55").Value = -684
$ws.Range("H113").Value = 6432.6665
$ws.Range("J113").Value = 6432.6665
$ws.Range("L113").Value = 6432.6665
$ws.Range("N113").Value = -12940.6665

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 882.375
$ws.Range("I2").Value = 777.5454999999999
$ws.Range("J2").Value = 1113
$ws.Range("K2").Value = 777.5454999999999
$ws.Range("L2").Value = 1113
$ws.Range("M2").Value = -664.5454999999999
$ws.Range("N2").Value = -1339
$ws.Range("H55").Value = 99999
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").Value = ""
$ws.Range("H116").Value = 882.375
$ws.Range("I116").Value = 777.5454999999999
$ws.Range("J116").Value = 1113
$ws.Range("K116").Value = 777.5454999999999
$ws.Range("L116").Value = 1113
$ws.Range("M116").Value = 1516.4545
$ws.Range("N116").Value = -5701
$ws.Range("H122").Value = 4742.25
$ws.Range("I122").Value = 4487.5
$ws.Range("K122").Value = 13462.5
$ws.Range("M122").Value = -11012.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 882.375
$ws.Range("I3").Value = 777.5454999999999
$ws.Range("J3").Value = 1113
$ws.Range("K3").Value = 777.5454999999999
$ws.Range("L3").Value = 1113
$ws.Range("M3").Value = -663.5454999999999
$ws.Range("N3").Value = -1341
$ws.Range("H20").Value = 1714.3334
$ws.Range("I20").Value = 1169.25
$ws.Range("J20").Value = 2804.5
$ws.Range("K20").Value = 1169.25
$ws.Range("L20").Value = 2804.5
$ws.Range("M20").Value = -922.25
$ws.Range("N20").Value = -3298.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 250
$ws.Range("I22").Value = 250
$ws.Range("K22").Value = 250
$ws.Range("M22").Value = 100
$ws.Range("H31").Value = 4619.6553
$ws.Range("I31").Value = 3187.7273
$ws.Range("K31").Value = 3187.7273
$ws.Range("M31").Value = -2892.7273
$ws.Range("H34").Value = 4619.6553
$ws.Range("I34").Value = 3187.7273
$ws.Range("K34").Value = 3187.7273
$ws.Range("M34").Value = -2985.7273
$ws.Range("H59").Value = 247002600
$ws.Range("J59").Value = 247002600
$ws.Range("L59").Value = 247002600
$ws.Range("N59").Value = -247004890
$ws.Range("H74").Value = 61882.5
$ws.Range("I74").Value = 25000
$ws.Range("J74").Value = 98765
$ws.Range("K74").Value = 25000
$ws.Range("L74").Value = 98765
$ws.Range("M74").Value = -24126
$ws.Range("N74").Value = -100513
$ws.Range("H77").Value = 61882.5
$ws.Range("I77").Value = 25000
$ws.Range("J77").Value = 98765
$ws.Range("K77").Value = 75000
$ws.Range("L77").Value = 296295
$ws.Range("M77").Value = -70632
$ws.Range("N77").Value = -305031
$ws.Range("H94").Value = 3962.9167
$ws.Range("I94").Value = 1924
$ws.Range("J94").Value = 8040.75
$ws.Range("K94").Value = 1924
$ws.Range("L94").Value = 8040.75
$ws.Range("M94").Value = -1473
$ws.Range("N94").Value = -8942.75
$ws.Range("H132").Value = 2572.2
$ws.Range("J132").Value = 3999.6667
$ws.Range("L132").Value = 11999.0001
$ws.Range("N132").Value = -17059.0001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1608.1666
$ws.Range("I51").Value = 1405
$ws.Range("K51").Value = 4215
$ws.Range("M51").Value = -3755
$ws.Range("H57").Value = 1500
$ws.Range("I57").Value = 1200
$ws.Range("J57").Value = 1650
$ws.Range("K57").Value = 3600
$ws.Range("L57").Value = 4950
$ws.Range("M57").Value = -3041
$ws.Range("N57").Value = -6068
$ws.Range("H59").Value = 356.66666
$ws.Range("I59").Value = 356.66666
$ws.Range("K59").Value = 1069.99998
$ws.Range("M59").Value = -529.9999800000001
$ws.Range("H60").Value = 1453.8235
$ws.Range("I60").Value = 87
$ws.Range("J60").Value = 2023.3334
$ws.Range("K60").Value = 261
$ws.Range("L60").Value = 6070.0002
$ws.Range("M60").Value = -10
$ws.Range("N60").Value = -6572.0002
$ws.Range("H80").Value = 4047.8147
$ws.Range("I80").Value = 3749.95
$ws.Range("J80").Value = 4898.857
$ws.Range("K80").Value = 11249.85
$ws.Range("L80").Value = 14696.571
$ws.Range("M80").Value = -10313.85
$ws.Range("N80").Value = -16568.571
$ws.Range("H83").Value = 4047.8147
$ws.Range("I83").Value = 3749.95
$ws.Range("J83").Value = 4898.857
$ws.Range("K83").Value = 33749.55
$ws.Range("L83").Value = 44089.713
$ws.Range("M83").Value = -29069.55
$ws.Range("N83").Value = -53449.713
$ws.Range("H107").Value = 693.6667
$ws.Range("I107").Value = 693.6667
$ws.Range("K107").Value = 2081.0001
$ws.Range("M107").Value = -161.0001000000002
$ws.Range("H112").Value = 913.25
$ws.Range("I112").Value = 913.25
$ws.Range("K112").Value = 2739.75
$ws.Range("M112").Value = -1631.75
$ws.Range("H131").Value = 2072.75
$ws.Range("I131").Value = 960.8570999999999
$ws.Range("J131").Value = 2530.5881
$ws.Range("K131").Value = 2882.5713
$ws.Range("L131").Value = 7591.7643
$ws.Range("M131").Value = 2157.4287
$ws.Range("N131").Value = -17671.7643
$ws.Range("H137").Value = 3516.5
$ws.Range("I137").Value = 3000
$ws.Range("K137").Value = 9000
$ws.Range("M137").Value = -3900

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 37559.11
$ws.Range("J29").Value = 34728.855
$ws.Range("L29").Value = 34728.855
$ws.Range("N29").Value = -35308.855
$ws.Range("H113").Value = 6044.7
$ws.Range("I113").Value = 3187.5
$ws.Range("J113").Value = 7949.5
$ws.Range("K113").Value = 3187.5
$ws.Range("L113").Value = 7949.5
$ws.Range("M113").Value = -1017.5
$ws.Range("N113").Value = -12289.5
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = ""

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 895.7
$ws.Range("I22").Value = 735.3333
$ws.Range("K22").Value = 735.3333
$ws.Range("M22").Value = -440.3333
$ws.Range("H27").Value = 895.7
$ws.Range("I27").Value = 735.3333
$ws.Range("K27").Value = 735.3333
$ws.Range("M27").Value = -628.3333
$ws.Range("H29").Value = 20000
$ws.Range("I29").Value = 20000
$ws.Range("K29").Value = 20000
$ws.Range("M29").Value = -19705
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").Value = ""
$ws.Range("H46").Value = 1666.6666
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").Value = ""
$ws.Range("H55").Value = 927.6842
$ws.Range("I55").Value = 897.7
$ws.Range("J55").Value = 961
$ws.Range("K55").Value = 897.7
$ws.Range("L55").Value = 961
$ws.Range("M55").Value = -724.7
$ws.Range("N55").Value = -1307
$ws.Range("H82").Value = 3374.5
$ws.Range("J82").Value = 999
$ws.Range("L82").Value = 999
$ws.Range("N82").Value = -1721
$ws.Range("H85").Value = 3374.5
$ws.Range("J85").Value = 999
$ws.Range("L85").Value = 999
$ws.Range("N85").Value = -3495
$ws.Range("H100").Value = 9166.700000000001
$ws.Range("I100").Value = 7916.75
$ws.Range("K100").Value = 7916.75
$ws.Range("M100").Value = -7375.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 69000
$ws.Range("J46").Value = 69000
$ws.Range("L46").Value = 69000
$ws.Range("N46").Value = -69462
$ws.Range("H54").Value = 99999
$ws.Range("J54").Value = 99999
$ws.Range("L54").Value = 99999
$ws.Range("N54").Value = -101039
$ws.Range("H81").Value = 3926.5557
$ws.Range("I81").Value = 3926.5557
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 7853.1114
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -6792.1114
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 3926.5557
$ws.Range("I84").Value = 3926.5557
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 39265.557
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -33961.557
$ws.Range("N84").Value = ""
$ws.Range("H122").Value = 3158.3635
$ws.Range("I122").Value = 3158.3635
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9475.0905
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7025.0905
$ws.Range("N122").Value = ""
$ws.Range("H126").Value = 5088.6665
$ws.Range("I126").Value = 2896.8
$ws.Range("J126").Value = 6654.2856
$ws.Range("K126").Value = 8690.400000000001
$ws.Range("L126").Value = 19962.8568
$ws.Range("M126").Value = -6220.400000000001
$ws.Range("N126").Value = -24902.8568
$ws.Range("H134").Value = 69000
$ws.Range("J134").Value = 69000
$ws.Range("L134").Value = 207000
$ws.Range("N134").Value = -212070
